# Natmi following Dr Hou advice
# A new "ECs" target cluster is introduced into the Gdnf-Ret NATMI output.
# This shifts the old row 2 (target=FAPs) stats to a recomputed row 3
# (target=FAPs) and the old row 3 (target=sCs) stats to a new row 4
# (target=sCs); row 2 becomes the new target=ECs row. All summary
# statistics are recalculated to reflect the extra cluster.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2: Sending cluster sCs -> Ligand Gdnf -> Receptor Ret -> Target ECs
$ws.Range("A2").Value = "sCs"
$ws.Range("B2").Value = "Gdnf"
$ws.Range("C2").Value = "Ret"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.538334
$ws.Range("H2").Value = 7.615002
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 2.027884333333333
$ws.Range("N2").Value = 6.083653
$ws.Range("O2").Value = 0.07019303724735296
$ws.Range("P2").Value = 0.07019303724735297
$ws.Range("Q2").Value = 5.147447751367332
$ws.Range("R2").Value = 46.327029762306
$ws.Range("S2").Value = 0.07019303724735296
$ws.Range("T2").Value = 0.07019303724735297

# ---- Row 3: Sending cluster sCs -> Ligand Gdnf -> Receptor Ret -> Target FAPs
$ws.Range("A3").Value = "sCs"
$ws.Range("B3").Value = "Gdnf"
$ws.Range("C3").Value = "Ret"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.538334
$ws.Range("H3").Value = 7.615002
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 25.63013966666666
$ws.Range("N3").Value = 76.89041899999999
$ws.Range("O3").Value = 0.887159745112283
$ws.Range("P3").Value = 0.887159745112283
$ws.Range("Q3").Value = 65.05785494064865
$ws.Range("R3").Value = 585.5206944658379
$ws.Range("S3").Value = 0.887159745112283
$ws.Range("T3").Value = 0.887159745112283

# ---- Row 4 (new): Sending cluster sCs -> Ligand Gdnf -> Receptor Ret -> Target sCs
$ws.Range("A4").Value = "sCs"
$ws.Range("B4").Value = "Gdnf"
$ws.Range("C4").Value = "Ret"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.538334
$ws.Range("H4").Value = 7.615002
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.232082666666667
$ws.Range("N4").Value = 3.696248
$ws.Range("O4").Value = 0.04264721764036409
$ws.Range("P4").Value = 0.04264721764036408
$ws.Range("Q4").Value = 3.127437323610666
$ws.Range("R4").Value = 28.146935912496
$ws.Range("S4").Value = 0.04264721764036409
$ws.Range("T4").Value = 0.04264721764036408
